$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt4"
$ws.Range("C2").Value = "Fzd2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.949238666666667
$ws.Range("H2").Value = 8.847716
$ws.Range("I2").Value = 0.5666145353461176
$ws.Range("J2").Value = 0.5666145353461176
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1278803333333333
$ws.Range("N2").Value = 0.383641
$ws.Range("O2").Value = 0.009974564977605908
$ws.Range("P2").Value = 0.009974564977605908
$ws.Range("Q2").Value = 0.3771496237728889
$ws.Range("R2").Value = 3.394346613956
$ws.Range("S2").Value = 0.005651733500065829
$ws.Range("T2").Value = 0.005651733500065829
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt4"
$ws.Range("C3").Value = "Fzd2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.949238666666667
$ws.Range("H3").Value = 8.847716
$ws.Range("I3").Value = 0.5666145353461176
$ws.Range("J3").Value = 0.5666145353461176
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.61985133333333
$ws.Range("N3").Value = 31.859554
$ws.Range("O3").Value = 0.8283400145723324
$ws.Range("P3").Value = 0.8283400145723324
$ws.Range("Q3").Value = 31.32047618651822
$ws.Range("R3").Value = 281.884285678664
$ws.Range("S3").Value = 0.4693494924654984
$ws.Range("T3").Value = 0.4693494924654984
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt4"
$ws.Range("C4").Value = "Fzd2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.949238666666667
$ws.Range("H4").Value = 8.847716
$ws.Range("I4").Value = 0.5666145353461176
$ws.Range("J4").Value = 0.5666145353461176
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.072911
$ws.Range("N4").Value = 6.218733
$ws.Range("O4").Value = 0.1616854204500617
$ws.Range("P4").Value = 0.1616854204500617
$ws.Range("Q4").Value = 6.113509273758666
$ws.Range("R4").Value = 55.021583463828
$ws.Range("S4").Value = 0.09161330938055336
$ws.Range("T4").Value = 0.09161330938055336
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt4"
$ws.Range("C5").Value = "Fzd2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.220831
$ws.Range("H5").Value = 3.662493
$ws.Range("I5").Value = 0.2345488676855596
$ws.Range("J5").Value = 0.2345488676855596
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1278803333333333
$ws.Range("N5").Value = 0.383641
$ws.Range("O5").Value = 0.009974564977605908
$ws.Range("P5").Value = 0.009974564977605908
$ws.Range("Q5").Value = 0.1561202752236667
$ws.Range("R5").Value = 1.405082477013
$ws.Range("S5").Value = 0.002339522921153505
$ws.Range("T5").Value = 0.002339522921153505
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt4"
$ws.Range("C6").Value = "Fzd2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.220831
$ws.Range("H6").Value = 3.662493
$ws.Range("I6").Value = 0.2345488676855596
$ws.Range("J6").Value = 0.2345488676855596
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.61985133333333
$ws.Range("N6").Value = 31.859554
$ws.Range("O6").Value = 0.8283400145723324
$ws.Range("P6").Value = 0.8283400145723324
$ws.Range("Q6").Value = 12.96504372312467
$ws.Range("R6").Value = 116.685393508122
$ws.Range("S6").Value = 0.1942862124765805
$ws.Range("T6").Value = 0.1942862124765805
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt4"
$ws.Range("C7").Value = "Fzd2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.220831
$ws.Range("H7").Value = 3.662493
$ws.Range("I7").Value = 0.2345488676855596
$ws.Range("J7").Value = 0.2345488676855596
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.072911
$ws.Range("N7").Value = 6.218733
$ws.Range("O7").Value = 0.1616854204500617
$ws.Range("P7").Value = 0.1616854204500617
$ws.Range("Q7").Value = 2.530674009041
$ws.Range("R7").Value = 22.776066081369
$ws.Range("S7").Value = 0.03792313228782559
$ws.Range("T7").Value = 0.03792313228782559
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Wnt4"
$ws.Range("C8").Value = "Fzd2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.034948
$ws.Range("H8").Value = 3.104844
$ws.Range("I8").Value = 0.1988365969683228
$ws.Range("J8").Value = 0.1988365969683228
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1278803333333333
$ws.Range("N8").Value = 0.383641
$ws.Range("O8").Value = 0.009974564977605908
$ws.Range("P8").Value = 0.009974564977605908
$ws.Range("Q8").Value = 0.1323494952226667
$ws.Range("R8").Value = 1.191145457004
$ws.Range("S8").Value = 0.001983308556386574
$ws.Range("T8").Value = 0.001983308556386574
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Wnt4"
$ws.Range("C9").Value = "Fzd2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.034948
$ws.Range("H9").Value = 3.104844
$ws.Range("I9").Value = 0.1988365969683228
$ws.Range("J9").Value = 0.1988365969683228
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.61985133333333
$ws.Range("N9").Value = 31.859554
$ws.Range("O9").Value = 0.8283400145723324
$ws.Range("P9").Value = 0.8283400145723324
$ws.Range("Q9").Value = 10.99099389773067
$ws.Range("R9").Value = 98.91894507957599
$ws.Range("S9").Value = 0.1647043096302535
$ws.Range("T9").Value = 0.1647043096302535
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Wnt4"
$ws.Range("C10").Value = "Fzd2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.034948
$ws.Range("H10").Value = 3.104844
$ws.Range("I10").Value = 0.1988365969683228
$ws.Range("J10").Value = 0.1988365969683228
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.072911
$ws.Range("N10").Value = 6.218733
$ws.Range("O10").Value = 0.1616854204500617
$ws.Range("P10").Value = 0.1616854204500617
$ws.Range("Q10").Value = 2.145355093628
$ws.Range("R10").Value = 19.308195842652
$ws.Range("S10").Value = 0.03214897878168273
$ws.Range("T10").Value = 0.03214897878168273
